# "Added trace files and finalized SA"
# Update the simulated-annealing (LS1-SA-RS) trace value for the UMissouri
# instance (row 12) with the finalized run's result, then let Excel's
# formula engine recompute the dependent %Over column and the summary
# Mean/Stdev statistics for column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E12").Value = 1220743

# Leave the cursor where it was left at save time.
$ws.Range("F2").Select()
